# Edit the two opening diary paragraphs: split their single run into
# several runs with revised wording, per the target diff.

$d = $word.ActiveDocument

function Get-ParagraphByPrefix {
    param([string]$prefix)
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# --- Paragraph 1: "This week we had to design a game ..." ---
$para1 = Get-ParagraphByPrefix("This week we had to design a game")
if ($para1 -eq $null) { throw "Paragraph 1 (This week...) not found" }
$range1 = $para1.Range
# Exclude the trailing paragraph mark so paragraph identity/formatting is kept.
$inner1 = $d.Range($range1.Start, $range1.End - 1)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This week </w:t></w:r><w:r><w:t>the group</w:t></w:r><w:r><w:t xml:space="preserve"> had to design a game which includes territorial acquisition, trading and </w:t></w:r><w:r><w:t>alliances (</w:t></w:r><w:r><w:t xml:space="preserve">team play). </w:t></w:r><w:r><w:t>The said group</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>was</w:t></w:r><w:r><w:t xml:space="preserve"> off campus and had to do the covid-19 session, which means the task that was given was different than the students who were present. </w:t></w:r><w:r><w:t>Online resources such as Wikipedia were used to get inspirations for the game.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$inner1.InsertXML($xml1)

# --- Paragraph 2: "The name of the game our team designed ..." ---
$para2 = Get-ParagraphByPrefix("The name of the game our team designed")
if ($para2 -eq $null) { throw "Paragraph 2 (The name of the game...) not found" }
$range2 = $para2.Range
$inner2 = $d.Range($range2.Start, $range2.End - 1)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The name of the game was Medieval Conquest. </w:t></w:r><w:r><w:t xml:space="preserve">The team were </w:t></w:r><w:r><w:t xml:space="preserve">asked to design a game which includes grabbing territory and team play, the team thought of the idea of a medieval type of game. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$inner2.InsertXML($xml2)

Write-Host "Paragraph 1 now: $($para1.Range.Text)"
Write-Host "Paragraph 2 now: $($para2.Range.Text)"
